$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044776057161285
$ws.Range("D2").Value = 1.052802739134675
$ws.Range("E2").Value = 1.058129198154092
$ws.Range("F2").Value = 1.065014483352743
$ws.Range("I2").Value = 1.043288728936267
$ws.Range("J2").Value = 1.049839333033683
$ws.Range("K2").Value = 1.055550463193599
$ws.Range("L2").Value = 1.060862276782518
$ws.Range("M2").Value = 1.067728865395869
$ws.Range("N2").Value = 1.020477565263184

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045622705547721
$ws.Range("D3").Value = 1.053480837374018
$ws.Range("E3").Value = 1.058939523776214
$ws.Range("F3").Value = 1.065852047513001
$ws.Range("I3").Value = 1.043483400266941
$ws.Range("J3").Value = 1.050333882869711
$ws.Range("K3").Value = 1.056041825576576
$ws.Range("L3").Value = 1.061486583960566
$ws.Range("M3").Value = 1.068381690695944
$ws.Range("N3").Value = 1.020643677651902

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046171274775451
$ws.Range("D4").Value = 1.053920248417721
$ws.Range("E4").Value = 1.059464928762126
$ws.Range("F4").Value = 1.066395116232679
$ws.Range("I4").Value = 1.043608534210297
$ws.Range("J4").Value = 1.050653942765955
$ws.Range("K4").Value = 1.056359723662849
$ws.Range("L4").Value = 1.061890966326332
$ws.Range("M4").Value = 1.068804569305417
$ws.Range("N4").Value = 1.020751134823777

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046402066638906
$ws.Range("D5").Value = 1.05410512730413
$ws.Range("E5").Value = 1.059686063294328
$ws.Range("F5").Value = 1.066623685777905
$ws.Range("I5").Value = 1.043660940757386
$ws.Range("J5").Value = 1.050788507367832
$ws.Range("K5").Value = 1.056493355593599
$ws.Range("L5").Value = 1.062061066193609
$ws.Range("M5").Value = 1.068982455236444
$ws.Range("N5").Value = 1.020796302476157

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046440827725802
$ws.Range("D6").Value = 1.054136178054469
$ws.Range("E6").Value = 1.059723207610566
$ws.Range("F6").Value = 1.066662079002191
$ws.Range("I6").Value = 1.043669728312234
$ws.Range("J6").Value = 1.050811101988973
$ws.Range("K6").Value = 1.056515792214759
$ws.Range("L6").Value = 1.062089632383065
$ws.Range("M6").Value = 1.069012329338639
$ws.Range("N6").Value = 1.020803885879101

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046174357948524
$ws.Range("D7").Value = 1.053922718189426
$ws.Range("E7").Value = 1.059467882575839
$ws.Range("F7").Value = 1.066398169357403
$ws.Range("I7").Value = 1.043609235254968
$ws.Range("J7").Value = 1.050655740781485
$ws.Range("K7").Value = 1.05636150930977
$ws.Range("L7").Value = 1.061893238827854
$ws.Range("M7").Value = 1.068806945804554
$ws.Range("N7").Value = 1.020751738385618

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045062033562394
$ws.Range("D8").Value = 1.053031772804865
$ws.Range("E8").Value = 1.058402828765259
$ws.Range("F8").Value = 1.0652973113336
$ws.Range("I8").Value = 1.043354690743699
$ws.Range("J8").Value = 1.050006456701278
$ws.Range("K8").Value = 1.055716530168193
$ws.Range("L8").Value = 1.061073177681434
$ws.Range("M8").Value = 1.067949395154072
$ws.Range("N8").Value = 1.020533709364449

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043107649212991
$ws.Range("D9").Value = 1.051466765298701
$ws.Range("E9").Value = 1.056534348408955
$ws.Range("F9").Value = 1.063366037342424
$ws.Range("I9").Value = 1.042899814899213
$ws.Range("J9").Value = 1.048862798685203
$ws.Range("K9").Value = 1.054579700470353
$ws.Range("L9").Value = 1.059631364070948
$ws.Range("M9").Value = 1.066441853892058
$ws.Range("N9").Value = 1.020149316317932

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041808640325835
$ws.Range("D10").Value = 1.050426864564136
$ws.Range("E10").Value = 1.055294376531469
$ws.Range("F10").Value = 1.06208440675669
$ws.Range("I10").Value = 1.042592347410211
$ws.Range("J10").Value = 1.048100748461507
$ws.Range("K10").Value = 1.053821697352177
$ws.Range("L10").Value = 1.058672425947386
$ws.Range("M10").Value = 1.065439326342039
$ws.Range("N10").Value = 1.01989294997852

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04124710445903
$ws.Range("D11").Value = 1.04997741428895
$ws.Range("E11").Value = 1.054758826818891
$ws.Range("F11").Value = 1.061530867398294
$ws.Range("I11").Value = 1.042458219459937
$ws.Range("J11").Value = 1.047770880657274
$ws.Range("K11").Value = 1.053493463085658
$ws.Range("L11").Value = 1.058257753504044
$ws.Range("M11").Value = 1.065005834508116
$ws.Range("N11").Value = 1.019781921875718

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041038668479695
$ws.Range("D12").Value = 1.049810595540146
$ws.Range("E12").Value = 1.054560107043426
$ws.Range("F12").Value = 1.061325472864701
$ws.Range("I12").Value = 1.042408249973556
$ws.Range("J12").Value = 1.047648370143836
$ws.Range("K12").Value = 1.053371541552631
$ws.Range("L12").Value = 1.058103810534324
$ws.Range("M12").Value = 1.064844909541036
$ws.Range("N12").Value = 1.019740678612931

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041083372186565
$ws.Range("D13").Value = 1.049846372946943
$ws.Range("E13").Value = 1.05460272370058
$ws.Range("F13").Value = 1.061369520944451
$ws.Range("I13").Value = 1.042418975301581
$ws.Range("J13").Value = 1.047674648271124
$ws.Range("K13").Value = 1.053397694145407
$ws.Range("L13").Value = 1.058136827959175
$ws.Range("M13").Value = 1.064879424248308
$ws.Range("N13").Value = 1.019749525537235

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041229872139979
$ws.Range("D14").Value = 1.049963622395707
$ws.Range("E14").Value = 1.054742396341416
$ws.Range("F14").Value = 1.061513885031989
$ws.Range("I14").Value = 1.042454091989628
$ws.Range("J14").Value = 1.047760753544108
$ws.Range("K14").Value = 1.05348338502025
$ws.Range("L14").Value = 1.05824502679024
$ws.Range("M14").Value = 1.064992530489176
$ws.Range("N14").Value = 1.019778512742542

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041320154654601
$ws.Range("D15").Value = 1.05003588054763
$ws.Range("E15").Value = 1.054828480784332
$ws.Range("F15").Value = 1.061602861012015
$ws.Range("I15").Value = 1.042475708898057
$ws.Range("J15").Value = 1.047813808146487
$ws.Range("K15").Value = 1.053536181943325
$ws.Range("L15").Value = 1.058311702950723
$ws.Range("M15").Value = 1.06506223137731
$ws.Range("N15").Value = 1.019796372402233

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04184592763533
$ws.Range("D16").Value = 1.050456710820182
$ws.Range("E16").Value = 1.055329948155894
$ws.Range("F16").Value = 1.062121173328831
$ws.Range("I16").Value = 1.04260122819097
$ws.Range("J16").Value = 1.048122643024336
$ws.Range("K16").Value = 1.05384348100394
$ws.Range("L16").Value = 1.058699958201806
$ws.Range("M16").Value = 1.065468108745658
$ws.Range("N16").Value = 1.019900318173066

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042175984940201
$ws.Range("D17").Value = 1.050720911021926
$ws.Range("E17").Value = 1.055644872570716
$ws.Range("F17").Value = 1.062446677252269
$ws.Range("I17").Value = 1.042679697913395
$ws.Range("J17").Value = 1.048316396009145
$ws.Range("K17").Value = 1.054036238833697
$ws.Range("L17").Value = 1.058943649933786
$ws.Range("M17").Value = 1.065722869180456
$ws.Range("N17").Value = 1.019965515658104

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042368592707966
$ws.Range("D18").Value = 1.050875094836139
$ws.Range("E18").Value = 1.055828694349595
$ws.Range("F18").Value = 1.062636674588914
$ws.Range("I18").Value = 1.042725372164882
$ws.Range("J18").Value = 1.048429418887751
$ws.Range("K18").Value = 1.054148669779567
$ws.Range("L18").Value = 1.059085844553117
$ws.Range("M18").Value = 1.065871525299581
$ws.Range("N18").Value = 1.020003542306104

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042434282350738
$ws.Range("D19").Value = 1.05092768109722
$ws.Range("E19").Value = 1.055891395103983
$ws.Range("F19").Value = 1.062701481845934
$ws.Range("I19").Value = 1.042740929631459
$ws.Range("J19").Value = 1.048467958431935
$ws.Range("K19").Value = 1.054187005541296
$ws.Range("L19").Value = 1.05913433824515
$ws.Range("M19").Value = 1.065922223094719
$ws.Range("N19").Value = 1.02001650806605

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042140563497961
$ws.Range("D20").Value = 1.05069255652079
$ws.Range("E20").Value = 1.055611070525347
$ws.Range("F20").Value = 1.062411739638217
$ws.Range("I20").Value = 1.042671288759712
$ws.Range("J20").Value = 1.048295607099572
$ws.Range("K20").Value = 1.054015557881294
$ws.Range("L20").Value = 1.058917498586973
$ws.Range("M20").Value = 1.06569552971648
$ws.Range("N20").Value = 1.01995852078231

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041186727581802
$ws.Range("D21").Value = 1.04992909183207
$ws.Range("E21").Value = 1.054701260493528
$ws.Range("F21").Value = 1.061471367465308
$ws.Range("I21").Value = 1.042443755090871
$ws.Range("J21").Value = 1.047735397205451
$ws.Range("K21").Value = 1.053458151201973
$ws.Range("L21").Value = 1.058213162576511
$ws.Range("M21").Value = 1.064959220934082
$ws.Range("N21").Value = 1.019769976797638

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040587843843738
$ws.Range("D22").Value = 1.049449807531433
$ws.Range("E22").Value = 1.054130426774923
$ws.Range("F22").Value = 1.060881360943118
$ws.Range("I22").Value = 1.042299837439967
$ws.Range("J22").Value = 1.047383270456184
$ws.Range("K22").Value = 1.053107684063878
$ws.Range("L22").Value = 1.057770810213065
$ws.Range("M22").Value = 1.064496814441434
$ws.Range("N22").Value = 1.019651417594574

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040905244134112
$ws.Range("D23").Value = 1.049703814800251
$ws.Range("E23").Value = 1.054432922131997
$ws.Range("F23").Value = 1.061194016087105
$ws.Range("I23").Value = 1.042376212017029
$ws.Range("J23").Value = 1.0475699296459
$ws.Range("K23").Value = 1.053293473189944
$ws.Range("L23").Value = 1.058005262436541
$ws.Range("M23").Value = 1.06474189312131
$ws.Range("N23").Value = 1.019714269245606

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042156568642272
$ws.Range("D24").Value = 1.050705368451457
$ws.Range("E24").Value = 1.055626343807994
$ws.Range("F24").Value = 1.062427526020929
$ws.Range("I24").Value = 1.04267508879042
$ws.Range("J24").Value = 1.048305000682101
$ws.Range("K24").Value = 1.054024902718093
$ws.Range("L24").Value = 1.058929315090272
$ws.Range("M24").Value = 1.065707883062379
$ws.Range("N24").Value = 1.019961681471654

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043612221848007
$ws.Range("D25").Value = 1.051870759261731
$ws.Range("E25").Value = 1.057016402302247
$ws.Range("F25").Value = 1.063864289237525
$ws.Range("I25").Value = 1.043018157625038
$ws.Range("J25").Value = 1.049158398856821
$ws.Range("K25").Value = 1.054873624694606
$ws.Range("L25").Value = 1.060003714011975
$ws.Range("M25").Value = 1.066831155985181
$ws.Range("N25").Value = 1.0202487116313

